$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete first data row (row 2); this shifts rows 3:53 up to 2:52,
# matching columns A, B and D for the refreshed forecast vectors.
$ws.Rows("2:2").Delete()

# Refresh the recomputed forecast columns (C = y_0_forecast, E = y_1_forecast)
# for every data row after the shift.
$cValues = @("", "", 1.514319819128396, 1.834695583582491, 2.007652128026982, 1.767835936772166, 1.022680528298392, 1.074400434091016, 0.9070039918702477, 0.9212998022035679, 1.022042907336096, 1.141837882844188, 1.34489417553354, 1.335361538769475, 1.277042522796856, 1.202048372526998, 2.615369162917314, 2.677488680362305, 2.618053282882693, 2.466954516646402, 1.731723847815725, 1.431088640641853, 1.372961566907027, 1.401189216021326, 1.983559881711905, 2.136062314641141, 2.241561867365394, 2.217567799050979, 2.041276490941102, 2.166968775134936, 2.139672475020404, 2.139672475020404, 2.330672672271739, 2.459440348120401, 2.100991693542231, 2.100991693542231, 0.7608230790701942, 0.8171929556848756, 0.8967077601845341, 0.8967077601845341, -0.4883557973630492, 0.8766015904249524, 0.782207885866093, 0.782207885866093, 1.617535832906758, 1.554086551645839, 1.508385007449875, 1.508385007449875, 0.3854686824285025, 0.5837948599211717, 0.6014263374495288)
$eValues = @(1.516248937663556, 1.560682679516057, 1.602279001294704, 2.1453644888767, 1.768040115052738, 1.317672174811868, 1.510468690286459, 1.501816644427989, 1.042579621507111, 1.028888107831327, 1.083482333436536, 1.303605130836716, 1.095916825800991, 1.192378712846454, 1.293136192195643, 1.210961441871872, 1.825134644920934, 2.033218171624651, 2.406099663413808, 2.152537330144288, 2.356276715023498, 2.21629047761287, 1.902399534782662, 1.966855307908655, 1.93103453922987, 2.031764787322499, 2.135927826705641, 1.950353221540246, 1.854752869950294, 1.984987808509886, 2.011395609719546, 2.210985773414453, 2.307457288603798, 2.526389380645511, 1.091147151778871, 1.114171399050901, 1.459415358104388, 1.509741350988136, 1.773412413757813, 0.1338254721205745, 0.4753196237801127, 2.477445663648559, 2.01691766737, 1.823564868738359, 1.142484412546874, 0.9888012784191602, 0.8780954941978392, 0.6266145540918089, 0.6767639290315763, 1.328924132093245, 1.609787824259601)

for ($i = 0; $i -lt $cValues.Length; $i++) {
    $row = $i + 2
    if ($cValues[$i] -ne "") {
        $ws.Cells.Item($row, 3).Value = $cValues[$i]
    } else {
        $ws.Cells.Item($row, 3).Value = ""
    }
    $ws.Cells.Item($row, 5).Value = $eValues[$i]
}
